$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Type" column header
$ws.Range("E1").Value = "Type"

# Domestic codes: rows 2 through 15 (DA01 .. DE01)
$ws.Range("E2:E15").Value = "Domestic"

# International codes: rows 16 through 18 (IF01, IG01, IG02)
$ws.Range("E16:E18").Value = "International"
